{"js": "// Fix three \"sprint to\" -> \"sprint two\" / \"manges\" -> \"manages\" typos.\nconst body = context.document.body;\n\nasync function replaceOnce(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replaceText, \"Replace\");\n    await context.sync();\n  }\n}\n\n// 1) \"During sprint to I worked\" -> \"During sprint two I worked\"\nawait replaceOnce(\"During sprint to I worked\", \"During sprint two I worked\");\n\n// 2) \"she manges to really help\" -> \"she manages to really help\"\nawait replaceOnce(\"she manges to really help\", \"she manages to really help\");\n\n// 3) \"3) Sprint to really made\" -> \"3) Sprint two really made\"\nawait replaceOnce(\"3) Sprint to really made\", \"3) Sprint two really made\");\n", "ps1": "# Fix three typos in the sprint two evaluation document:\n#   \"sprint to\"   -> \"sprint two\"   (x2)\n#   \"manges\"      -> \"manages\"\n$d = $word.ActiveDocument\n\n$r1 = $d.Content\n$r1.Find.Execute(\"During sprint to I worked\", $true, $false, $false, $false, $false, $true, 1, $false, \"During sprint two I worked\", 2)\n\n$r2 = $d.Content\n$r2.Find.Execute(\"she manges to really help\", $true, $false, $false, $false, $false, $true, 1, $false, \"she manages to really help\", 2)\n\n$r3 = $d.Content\n$r3.Find.Execute(\"3) Sprint to really made\", $true, $false, $false, $false, $false, $true, 1, $false, \"3) Sprint two really made\", 2)\n"}
